# origin-verification.md: updated appveyor instructions and diagram
#
# 1) Footer "updated" date field (06.02.2023 -> 03.03.2023) on the slide
#    master and every slide layout.
# 2) Slide 9 diagram rework: circled numerals instead of "(n)", labels
#    pushed further down the diagram, connectors re-routed, and the two
#    outline rectangles grown taller with a heavier border.

$EMU_PER_PT = 12700

$p = $ppt.ActivePresentation

function Set-DateFieldText($shapes, [string]$newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "*Date*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- 1) date field everywhere it is defined -------------------------------
Set-DateFieldText $p.SlideMaster.Shapes "03.03.2023"

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DateFieldText $layouts.Item($li).Shapes "03.03.2023"
}

# --- 2) slide 9 AppVeyor diagram -------------------------------------------
$s = $p.Slides.Item(9)

# (1) copy bearer token  -- numeral only, keep the rest of the label
$tb14 = $s.Shapes.Item(6)
$tb14.TextFrame.TextRange.Runs(1).Text = [char]0x2460

# "CI User" box moves down and gains "token" in its label
$tb15 = $s.Shapes.Item(7)
$tb15.Top = 2712643 / $EMU_PER_PT
$tb15.TextFrame.TextRange.Runs(1).Text = "CI User token"

# "Settings -> Encrypt YAML" box moves down (text unchanged)
$tb16 = $s.Shapes.Item(8)
$tb16.Top = 2712643 / $EMU_PER_PT

# connector from Bearer-token label down to CI User / Settings row
$conn17 = $s.Shapes.Item(9)
$conn17.Top = 2897309 / $EMU_PER_PT

# (3) copy encrypted value -- numeral only, moves down
$tb21 = $s.Shapes.Item(10)
$tb21.Top = 3312585 / $EMU_PER_PT
$tb21.TextFrame.TextRange.Runs(1).Text = [char]0x2462 + "`t"
# the autosize engine recomputes a shorter box after the run edit above;
# put the (unchanged) height back to its original value
$tb21.Height = 523220 / $EMU_PER_PT

# "appveyor.yml" label moves down
$tb22 = $s.Shapes.Item(11)
$tb22.Top = 3941186 / $EMU_PER_PT

# connector down into appveyor.yml moves down
$conn23 = $s.Shapes.Item(12)
$conn23.Top = 3312585 / $EMU_PER_PT

# (2) Bearer + API token -- numeral + trailing space merged into one run,
# box moves down
$tb29 = $s.Shapes.Item(13)
$tb29.Top = 2544239 / $EMU_PER_PT
$tb29.TextFrame.TextRange.Runs(2).Text = ""
$tb29.TextFrame.TextRange.Runs(1).Text = [char]0x2461 + " "

# left outline rectangle grows taller and gets a heavier border
$rect18 = $s.Shapes.Item(14)
$rect18.Height = 1408678 / $EMU_PER_PT
$rect18.Line.Weight = 1.5

# right outline rectangle grows taller and gets a heavier border
$rect24 = $s.Shapes.Item(15)
$rect24.Height = 1408678 / $EMU_PER_PT
$rect24.Line.Weight = 1.5

# git icon moves down
$pic2 = $s.Shapes.Item(16)
$pic2.Top = 3947385 / $EMU_PER_PT

# "secure: gYX+..." caption moves down
$tb27 = $s.Shapes.Item(17)
$tb27.Top = 4310641 / $EMU_PER_PT
